# The exported report's rows got re-sorted by "ID Venda" (column A) in
# ascending order. Re-apply that sort to the data rows (row 1 is the header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162
$xlToLeft = -4159
$xlAscending = 1

$headerRow = 1
$firstDataRow = $headerRow + 1
$lastRow = $ws.Cells($ws.Rows.Count, 1).End($xlUp).Row
$lastCol = $ws.Cells($headerRow, $ws.Columns.Count).End($xlToLeft).Column

$dataRange = $ws.Range($ws.Cells($firstDataRow, 1), $ws.Cells($lastRow, $lastCol))
$keyRange = $ws.Range($ws.Cells($firstDataRow, 1), $ws.Cells($lastRow, 1))

$dataRange.Sort($keyRange, $xlAscending)
